# Adding Horeco to the Portfolio
# Shift the date/time stamps in column A forward by 12 days (one full cycle
# of 96 quarter-hourly readings later), and update the production values in
# column B for the rows affected by the newly added Horeco plant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: every timestamp in rows 2-97 moves forward by 12 days.
# Use Value2 to read/write the raw serial date number (avoids the
# formatted-date/string round-trip that Value would otherwise perform).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 12
}

# Column B: updated production figures (rows 22-39) reflecting the new
# portfolio composition after adding Horeco.
$newB = @{
    22 = 0
    23 = 1
    24 = 18
    25 = 36
    26 = 74
    27 = 127
    28 = 198
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
